$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-6 from 45185 to 45204
# (keep existing date formatting / style untouched)
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45204
}
